$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 391; this shifts the existing rows 391-412 down to 392-413
# and carries the row-391 formatting (e.g. the date style on column D) down with it.
$ws.Rows("391:391").Insert()

# Populate the newly inserted row 391 with the new record.
$ws.Range("A391").Value = 5
$ws.Range("B391").Value = "Macroferia Regional de Talca"
$ws.Range("C391").Value = "Maule"
$ws.Range("D391").Value = 44516
$ws.Range("E391").Value = 7
$ws.Range("F391").Value = 100112004
$ws.Range("G391").Value = "Cebolla"
$ws.Range("H391").Value = "Sin especificar"
$ws.Range("I391").Value = "1a nueva(o)"
$ws.Range("J391").Value = 50000
$ws.Range("K391").Value = 1000
$ws.Range("L391").Value = 1000
$ws.Range("M391").Value = 1000
$ws.Range("N391").Value = "$/paquete 10 unidades (volumen en unidades)"
$ws.Range("O391").Value = "Región de O'Higgins"
$ws.Range("P391").Value = 100
$ws.Range("Q391").Value = 10
$ws.Range("R391").Value = "Hortaliza"
